$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the two rows that were folded into neighbouring rows:
#    - old row 16 "Partner Certificate Signing and RE-issueing" (merged into
#      the "Partner Certificate Validation" row)
#    - old row 11 "Policy ID Validation" (merged into the new
#      "Map Policies to Partners" row)
#    Delete from the bottom up so the earlier row number stays valid.
# ---------------------------------------------------------------------------
$ws.Rows(16).Delete()
$ws.Rows(11).Delete()

# ---------------------------------------------------------------------------
# 2. Fix up row heights for the rows whose wrapped text changed length.
#    (Rows below shifted up but kept their old height since Excel doesn't
#    auto-reflow heights on a plain delete.)
# ---------------------------------------------------------------------------
$ws.Rows(10).RowHeight = 145
$ws.Rows(11).RowHeight = 29
$ws.Rows(14).RowHeight = 145

# ---------------------------------------------------------------------------
# 3. Text corrections on rows that kept their place.
# ---------------------------------------------------------------------------

# Row 5: MISP License Key Pattern Validation -> VALID/INVALID wording
$ws.Range("E5").Value2 = "1. Validate length of a License Key as configured and respond as mentioned below`na. If found valid, respond with ""VALID""`nb. if found invalid, respond with ""INVALID"""

# Row 9: Partner ID Validation -> VALID/INVALID wording
$ws.Range("E9").Value2 = "1. Validate length of a Partner ID as configured and respond as mentioned below`na. If found valid, respond with ""VALID""`nb. if found invalid, respond with ""INVALID"""

# Row 10 (was "Policy ID Generation"): becomes "Map Policies to Partners"
$ws.Range("D10").Value2 = "Map Policies to Partners"
$ws.Range("E10").Value2 = "1. Map following Policies to Partners`na. Auth Policies ( can be Mandatory/Non-Mandatory)`n     1. OTP Trigger `n     2. OTP Authentication`n     3. Demo Authentication `n     4. Biometric Authentication - FMR Data Match `n     5. Biometric Authentication - IIR Data Match  `n     6. Biometric Authentication - FID Data Match `nb. E-Kyc Policies (can be Required/Not Required)`n    1. eKYC - all combinations of eKYC demo fields "

# Row 11 (was "Policy ID"): becomes "Retrieve Policies based on Partner ID"
$ws.Range("D11").Value2 = "Retrieve Policies based on Partner ID"
$ws.Range("E11").Value2 = "1. Receive request to retreive policies based on Partner ID`n2. Respond appropirately if Partner ID does not exist"

# Row 12: Partner Registration -> "Map Policies to the Partner" wording
$ws.Range("E12").Value2 = "1. Receive request to register a Partner with follwing parameters`na. Partner Name`nb. Partner Contact Name`nc. Partner Phone`nd. Partner Email ID`n2. Issue and Map Partner ID`n3. Map Policies to the Partner`na. Multiple Policies can be mapped to a Partner`nb. A Partner can have a policy for both Auth and E-KYC`n4. Store the Partner in MOSIP"

# Row 13: MISP - Partner Mapping -> fix "ve" -> "be" typo
$ws.Range("E13").Value2 = "1. Receive a request to map MISP to a Partner with MISP ID and Partner ID as Input`n2. There can be a many-to-mapping between MISPs and Partners"

# Row 14 (was "Partner Certiicate Validation"): fix typo + merge in the
# certificate signing/re-issuing content
$ws.Range("D14").Value2 = "Partner Certificate Validation"
$ws.Range("E14").Value2 = "1. Upload Digital Certificate on Admin Portal for a Partner`n2. Verify CA Authority of the certificate`n3. Sign the certificate with MOSIP Certificate`n4. Respond to the source with the re-issued certificate`n5. Certificate will be uploaded by the MOSIP admin. The Registered Partner will send the certificate to the MOSIP Admin through ofline process. Re-issued certificate will be sent to the Partner by MOSIP admin through notification/offline process`n6. Private key to change priodically as per the Key Rotation Policy set by admin"

# Row 15: Distribution of Public Key to Partners -> reworded
$ws.Range("E15").Value2 = "1. Distribute Public Key to Partners for encrypting the Auth Request befoe sending it to the MOSIP`n2. Public key needs to be distributed priodically whenever the correspinding Private Key is rotated"
